$d = $word.ActiveDocument
$wdFindContinue = 1
$wdReplaceAll = 2
$wdAlignParagraphJustify = 3

function Find-Replace($findText, $replaceText, $matchCase=$true, $matchWhole=$false) {
    $rng = $d.Content
    $rng.Find.Execute($findText, $matchCase, $matchWhole, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll) | Out-Null
}

# 1. Remove the stray _GoBack bookmark after "Versión 1.5" (Word renumbers
#    the remaining bookmark ids sequentially on save, which also produces
#    the _Toc bookmark id shifts seen in the diff).
$d.Bookmarks("_GoBack").Delete()

# 3. Table cell wording update (version history row).
Find-Replace "Verificación Ortográfica " "Corrección ortográfica luego de la realización de los ítems anteriores"

function Justify-ParagraphContaining($anchorText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
    if ($found) {
        $rng.Paragraphs(1).Format.Alignment = $wdAlignParagraphJustify
    }
    return $found
}

# 8-13. Justify the tool-comparison table description cells.
Justify-ParagraphContaining "Git es un sistema de control de versiones distribuido de código abierto" | Out-Null
Justify-ParagraphContaining "Es la herramienta que se utiliza" | Out-Null
Justify-ParagraphContaining "Permite crear los proyectos con JAVA primordialmente" | Out-Null
Justify-ParagraphContaining "YII es el framework " | Out-Null
Justify-ParagraphContaining "Es un SGBD " | Out-Null
Justify-ParagraphContaining "Es un servidor web de código abierto usado para desarrollar de forma " | Out-Null

# 9b. Add a trailing "." run/sentence after the Github-tool description cell.
$rng = $d.Content
$rng.Find.Execute("para alojar los proyectos utilizando el sistema de control de versiones GIT", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(".")

# 14. Remove the first of the (now) three empty Body Text paragraphs
#     that trail the tools-comparison table.
$toolsTable = $d.Tables(3)
$afterTableRng = $d.Range($toolsTable.Range.End, $toolsTable.Range.End)
$afterTableRng.Paragraphs(1).Range.Delete()

Write-Output "stage B done"
